$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("W12")

# Ryan (row 2): Pts and Proj both become 112.78
$ws.Range("B2").Value = 112.78
$ws.Range("C2").Value = 112.78

# James (row 9): Pts and Proj both become 104.84
$ws.Range("B9").Value = 104.84
$ws.Range("C9").Value = 104.84

# Mike (row 11): Pts and Proj both become 147.74
$ws.Range("B11").Value = 147.74
$ws.Range("C11").Value = 147.74

# Update the active selection to D10, matching the authored edit
$ws.Range("D10").Select()
